$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.2736413333333333
$ws.Range("H2").Value = 0.820924
$ws.Range("I2").Value = 0.1876387849846732
$ws.Range("J2").Value = 0.1876387849846732
$ws.Range("M2").Value = 2.866432
$ws.Range("N2").Value = 8.599295999999999
$ws.Range("O2").Value = 0.9456981836489474
$ws.Range("P2").Value = 0.9456981836489475
$ws.Range("Q2").Value = 0.7843742743893333
$ws.Range("R2").Value = 7.059368469503999
$ws.Range("S2").Value = 0.1774496581421008
$ws.Range("T2").Value = 0.1774496581421009

# Row 3
$ws.Range("G3").Value = 0.2736413333333333
$ws.Range("H3").Value = 0.820924
$ws.Range("I3").Value = 0.1876387849846732
$ws.Range("J3").Value = 0.1876387849846732
$ws.Range("O3").Value = 0.05430181635105255
$ws.Range("P3").Value = 0.05430181635105256
$ws.Range("Q3").Value = 0.04503862705333333
$ws.Range("R3").Value = 0.40534764348
$ws.Range("S3").Value = 0.01018912684257236
$ws.Range("T3").Value = 0.01018912684257236

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "ECs"
$ws.Range("G4").Value = 0.9293790000000001
$ws.Range("H4").Value = 2.788137
$ws.Range("I4").Value = 0.6372851068440097
$ws.Range("J4").Value = 0.6372851068440097
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 2.866432
$ws.Range("N4").Value = 8.599295999999999
$ws.Range("O4").Value = 0.9456981836489474
$ws.Range("P4").Value = 0.9456981836489475
$ws.Range("Q4").Value = 2.664001705728
$ws.Range("R4").Value = 23.976015351552
$ws.Range("S4").Value = 0.6026793680089054
$ws.Range("T4").Value = 0.6026793680089054

# Row 5
$ws.Range("D5").Value = "FAPs"
$ws.Range("I5").Value = 0.6372851068440097
$ws.Range("J5").Value = 0.6372851068440097
$ws.Range("K5").Value = 1.0
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.16459
$ws.Range("N5").Value = 0.49377
$ws.Range("O5").Value = 0.05430181635105255
$ws.Range("P5").Value = 0.05430181635105256
$ws.Range("Q5").Value = 0.15296648961
$ws.Range("R5").Value = 1.37669840649
$ws.Range("S5").Value = 0.03460573883510432
$ws.Range("T5").Value = 0.03460573883510432

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("D6").Value = "ECs"
$ws.Range("G6").Value = 0.2553206666666667
$ws.Range("H6").Value = 0.765962
$ws.Range("I6").Value = 0.1750761081713171
$ws.Range("J6").Value = 0.175076108171317
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 2.866432
$ws.Range("N6").Value = 8.599295999999999
$ws.Range("O6").Value = 0.9456981836489474
$ws.Range("P6").Value = 0.9456981836489475
$ws.Range("Q6").Value = 0.7318593291946667
$ws.Range("R6").Value = 6.586733962752
$ws.Range("S6").Value = 0.1655691574979412
$ws.Range("T6").Value = 0.1655691574979412

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("D7").Value = "FAPs"
$ws.Range("G7").Value = 0.2553206666666667
$ws.Range("H7").Value = 0.765962
$ws.Range("I7").Value = 0.1750761081713171
$ws.Range("J7").Value = 0.175076108171317
$ws.Range("K7").Value = 1.0
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.16459
$ws.Range("N7").Value = 0.49377
$ws.Range("O7").Value = 0.05430181635105255
$ws.Range("P7").Value = 0.05430181635105256
$ws.Range("Q7").Value = 0.04202322852666667
$ws.Range("R7").Value = 0.37820905674
$ws.Range("S7").Value = 0.009506950673375869
$ws.Range("T7").Value = 0.009506950673375869

# Drop the MuSCs-target rows (old rows 8-10); row count shrinks to match new dimension A1:T7
$ws.Rows("8:10").Delete()
